$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUSY-RF2")

# --- Apply the "2 decimal" number format (style s="1") to the MEDIA row 80 ---
$ws.Range("B80:V80").NumberFormat = "0.00"

# --- New block below the existing table: "Tiempos medios" ---

# Row 96: title (new shared string "Tiempos medios")
$ws.Range("B96").Value = "Tiempos medios"

# Row 97: header numbers (depth 4..10)
$ws.Range("B97").Value = 4
$ws.Range("E97").Value = 5
$ws.Range("H97").Value = 6
$ws.Range("K97").Value = 7
$ws.Range("N97").Value = 8
$ws.Range("Q97").Value = 9
$ws.Range("T97").Value = 10

# Row 98: averages of the ten MEDIA rows, formatted like row 80 (style s="1")
$ws.Range("B98:U98").NumberFormat = "0.00"

$ws.Range("B98").Formula = "=(B80+B72+B64+B56+B48+B40+B32+B24+B16+B8)/10"
$ws.Range("C98").Formula = "=(C80+C72+C64+C56+C48+C40+C32+C24+C16+C8)/10"
$ws.Range("E98:F98").Formula = "=(E80+E72+E64+E56+E48+E40+E32+E24+E16+E8)/10"
$ws.Range("H98:I98").Formula = "=(H80+H72+H64+H56+H48+H40+H32+H24+H16+H8)/10"
$ws.Range("K98:L98").Formula = "=(K80+K72+K64+K56+K48+K40+K32+K24+K16+K8)/10"
$ws.Range("N98:O98").Formula = "=(N80+N72+N64+N56+N48+N40+N32+N24+N16+N8)/10"
$ws.Range("Q98:R98").Formula = "=(Q80+Q72+Q64+Q56+Q48+Q40+Q32+Q24+Q16+Q8)/10"
$ws.Range("T98:U98").Formula = "=(T80+T72+T64+T56+T48+T40+T32+T24+T16+T8)/10"

# Row 99: CPU / GPU labels under each pair of columns
$ws.Range("B99").Value = "CPU"
$ws.Range("C99").Value = "GPU"
$ws.Range("E99").Value = "CPU"
$ws.Range("F99").Value = "GPU"
$ws.Range("H99").Value = "CPU"
$ws.Range("I99").Value = "GPU"
$ws.Range("K99").Value = "CPU"
$ws.Range("L99").Value = "GPU"
$ws.Range("N99").Value = "CPU"
$ws.Range("O99").Value = "GPU"
$ws.Range("Q99").Value = "CPU"
$ws.Range("R99").Value = "GPU"
$ws.Range("T99").Value = "CPU"
$ws.Range("U99").Value = "GPU"

# --- View-state bookkeeping matching the saved sheet view ---
[void]$ws.Range("B94").Select()
